$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = 0.2200772200772201
$ws.Cells.Item(2, 3).Value = 0.5096525096525096
$ws.Cells.Item(2, 10).Value = 0.0193050193050193
$ws.Cells.Item(2, 15).Value = 0.003861003861003861
$ws.Cells.Item(2, 16).Value = 0.1467181467181467
$ws.Cells.Item(2, 19).Value = 0.1003861003861004
$ws.Cells.Item(3, 2).Value = 0.02205882352941177
$ws.Cells.Item(3, 3).Value = 0.02205882352941177
$ws.Cells.Item(3, 10).Value = 0.03676470588235294
$ws.Cells.Item(3, 16).Value = 0.7279411764705882
$ws.Cells.Item(3, 19).Value = 0.1911764705882353
$ws.Cells.Item(4, 10).Value = 0.08163265306122448
$ws.Cells.Item(4, 16).Value = 0.6122448979591837
$ws.Cells.Item(4, 19).Value = 0.3061224489795918
$ws.Cells.Item(6, 2).Value = 0.05376344086021505
$ws.Cells.Item(6, 4).Value = 0.007168458781362007
$ws.Cells.Item(6, 5).Value = 0.003584229390681004
$ws.Cells.Item(6, 10).Value = 0.2007168458781362
$ws.Cells.Item(6, 15).Value = 0.01433691756272401
$ws.Cells.Item(6, 17).Value = 0.1362007168458781
$ws.Cells.Item(6, 18).Value = 0.06810035842293907
$ws.Cells.Item(6, 19).Value = 0.4050179211469534
$ws.Cells.Item(7, 2).Value = 0.1146788990825688
$ws.Cells.Item(7, 4).Value = 0.004587155963302753
$ws.Cells.Item(7, 6).Value = 0.0871559633027523
$ws.Cells.Item(7, 10).Value = 0.1284403669724771
$ws.Cells.Item(7, 15).Value = 0.009174311926605505
$ws.Cells.Item(7, 17).Value = 0.1697247706422018
$ws.Cells.Item(7, 18).Value = 0.08256880733944955
$ws.Cells.Item(7, 19).Value = 0.4036697247706422
$ws.Cells.Item(8, 2).Value = 0.08057851239669421
$ws.Cells.Item(8, 4).Value = 0.01033057851239669
$ws.Cells.Item(8, 6).Value = 0.06611570247933884
$ws.Cells.Item(8, 10).Value = 0.1239669421487603
$ws.Cells.Item(8, 15).Value = 0.006198347107438017
$ws.Cells.Item(8, 17).Value = 0.1735537190082645
$ws.Cells.Item(8, 18).Value = 0.07231404958677685
$ws.Cells.Item(8, 19).Value = 0.4669421487603306
$ws.Cells.Item(9, 2).Value = 0.0851063829787234
$ws.Cells.Item(9, 4).Value = 0.005319148936170213
$ws.Cells.Item(9, 5).Value = 0.005319148936170213
$ws.Cells.Item(9, 6).Value = 0.06382978723404255
$ws.Cells.Item(9, 10).Value = 0.1382978723404255
$ws.Cells.Item(9, 15).Value = 0.01595744680851064
$ws.Cells.Item(9, 17).Value = 0.1436170212765958
$ws.Cells.Item(9, 18).Value = 0.07446808510638298
$ws.Cells.Item(9, 19).Value = 0.4680851063829787
$ws.Cells.Item(10, 2).Value = 0.08333333333333333
$ws.Cells.Item(10, 4).Value = 0.03465346534653466
$ws.Cells.Item(10, 6).Value = 0.09075907590759076
$ws.Cells.Item(10, 10).Value = 0.1221122112211221
$ws.Cells.Item(10, 15).Value = 0.01897689768976898
$ws.Cells.Item(10, 17).Value = 0.2120462046204621
$ws.Cells.Item(10, 18).Value = 0.07508250825082509
$ws.Cells.Item(10, 19).Value = 0.363036303630363
$ws.Cells.Item(11, 7).Value = 0.1405228758169935
$ws.Cells.Item(11, 10).Value = 0.07843137254901961
$ws.Cells.Item(11, 11).Value = 0.1928104575163399
$ws.Cells.Item(11, 12).Value = 0.565359477124183
$ws.Cells.Item(11, 19).Value = 0.02287581699346405
$ws.Cells.Item(12, 7).Value = 0.76
$ws.Cells.Item(12, 10).Value = 0.1942857142857143
$ws.Cells.Item(12, 11).Value = 0.005714285714285714
$ws.Cells.Item(12, 12).Value = 0.01714285714285714
$ws.Cells.Item(12, 19).Value = 0.02285714285714286
$ws.Cells.Item(13, 7).Value = 0.7619047619047619
$ws.Cells.Item(13, 10).Value = 0.2063492063492063
$ws.Cells.Item(13, 19).Value = 0.03174603174603174
$ws.Cells.Item(15, 6).Value = 0.004219409282700422
$ws.Cells.Item(15, 8).Value = 0.1983122362869198
$ws.Cells.Item(15, 9).Value = 0.0759493670886076
$ws.Cells.Item(15, 10).Value = 0.3586497890295359
$ws.Cells.Item(15, 11).Value = 0.07172995780590717
$ws.Cells.Item(15, 13).Value = 0.01265822784810127
$ws.Cells.Item(15, 15).Value = 0.07172995780590717
$ws.Cells.Item(15, 19).Value = 0.2067510548523207
$ws.Cells.Item(16, 6).Value = 0.03067484662576687
$ws.Cells.Item(16, 8).Value = 0.2269938650306748
$ws.Cells.Item(16, 9).Value = 0.03680981595092025
$ws.Cells.Item(16, 10).Value = 0.3558282208588957
$ws.Cells.Item(16, 11).Value = 0.1165644171779141
$ws.Cells.Item(16, 13).Value = 0.01840490797546012
$ws.Cells.Item(16, 15).Value = 0.05521472392638037
$ws.Cells.Item(16, 19).Value = 0.1595092024539877
$ws.Cells.Item(17, 6).Value = 0.01348314606741573
$ws.Cells.Item(17, 8).Value = 0.2022471910112359
$ws.Cells.Item(17, 9).Value = 0.07191011235955057
$ws.Cells.Item(17, 10).Value = 0.3910112359550562
$ws.Cells.Item(17, 11).Value = 0.09662921348314607
$ws.Cells.Item(17, 13).Value = 0.01797752808988764
$ws.Cells.Item(17, 15).Value = 0.07191011235955057
$ws.Cells.Item(17, 19).Value = 0.1348314606741573
$ws.Cells.Item(18, 6).Value = 0.02857142857142857
$ws.Cells.Item(18, 8).Value = 0.1542857142857143
$ws.Cells.Item(18, 9).Value = 0.09714285714285714
$ws.Cells.Item(18, 10).Value = 0.36
$ws.Cells.Item(18, 11).Value = 0.1085714285714286
$ws.Cells.Item(18, 13).Value = 0.05142857142857143
$ws.Cells.Item(18, 14).Value = 0.005714285714285714
$ws.Cells.Item(18, 15).Value = 0.09714285714285714
$ws.Cells.Item(18, 19).Value = 0.09714285714285714
$ws.Cells.Item(19, 6).Value = 0.01779095626389918
$ws.Cells.Item(19, 8).Value = 0.2164566345441067
$ws.Cells.Item(19, 9).Value = 0.08673091178650852
$ws.Cells.Item(19, 10).Value = 0.3335804299481097
$ws.Cells.Item(19, 11).Value = 0.1074870274277242
$ws.Cells.Item(19, 13).Value = 0.03187546330615271
$ws.Cells.Item(19, 14).Value = 0.002223869532987398
$ws.Cells.Item(19, 15).Value = 0.07116382505559674
$ws.Cells.Item(19, 19).Value = 0.1326908821349148
